$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.1
$ws.Range("O3").Value = 1.54
$ws.Range("P3").Value = 2.25
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 1.1
$ws.Range("Y3").Value = 1.67
$ws.Range("M4").Value = 1.11
$ws.Range("O4").Value = 1.63
$ws.Range("X4").Value = 1.07
$ws.Range("Y4").Value = 1.73
$ws.Range("M5").Value = 1.16
$ws.Range("N5").Value = 4.6
$ws.Range("O5").Value = 1.7
$ws.Range("P5").Value = 2.02
$ws.Range("W5").Value = 5.8
$ws.Range("X5").Value = 1.1
$ws.Range("Y5").Value = 1.75
$ws.Range("Z5").Value = 2.05
$ws.Range("AA5").Value = 2.77
$ws.Range("AB5").Value = 1.39
$ws.Range("G6").Value = 2.5
$ws.Range("I6").Value = 2.75
$ws.Range("J6").Value = 3.1
$ws.Range("L6").Value = 3.25
$ws.Range("M6").Value = 1.03
$ws.Range("N6").Value = 11
$ws.Range("O6").Value = 1.22
$ws.Range("S6").Value = 1.85
$ws.Range("T6").Value = 2
$ws.Range("W6").Value = 3.2
$ws.Range("X6").Value = 1.33
$ws.Range("Y6").Value = 1.36
$ws.Range("AP6").Value = 10
$ws.Range("AS6").Value = 26
$ws.Range("G7").Value = 1.65
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 5.6
$ws.Range("J7").Value = 2.2
$ws.Range("K7").Value = 2.1
$ws.Range("M7").Value = 1.08
$ws.Range("N7").Value = 6.4
$ws.Range("O7").Value = 1.38
$ws.Range("P7").Value = 2.82
$ws.Range("S7").Value = 2.12
$ws.Range("T7").Value = 1.65
$ws.Range("W7").Value = 3.65
$ws.Range("AA7").Value = 2.05
$ws.Range("AB7").Value = 1.7
$ws.Range("AC7").Value = 5.5
$ws.Range("AD7").Value = 6.8
$ws.Range("AF7").Value = 12
$ws.Range("AI7").Value = 6.4
$ws.Range("AJ7").Value = 6.6
$ws.Range("AK7").Value = 18
$ws.Range("AN7").Value = 13
$ws.Range("AO7").Value = 32
$ws.Range("AP7").Value = 17.5
$ws.Range("AQ7").Value = 120
$ws.Range("AS7").Value = 65
$ws.Range("G9").Value = 1.39
$ws.Range("H9").Value = 4
$ws.Range("J9").Value = 1.87
$ws.Range("K9").Value = 2.25
$ws.Range("L9").Value = 7.7
$ws.Range("N9").Value = 6.9
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.05
$ws.Range("S9").Value = 1.98
$ws.Range("T9").Value = 1.75
$ws.Range("W9").Value = 3.3
$ws.Range("X9").Value = 1.28
$ws.Range("Y9").Value = 1.39
$ws.Range("Z9").Value = 2.77
$ws.Range("AA9").Value = 2.25
$ws.Range("AB9").Value = 1.57
$ws.Range("I10").Value = 3.7
$ws.Range("L10").Value = 4
$ws.Range("M10").Value = 1.02
$ws.Range("O10").Value = 1.13
$ws.Range("W10").Value = 2.37
$ws.Range("X10").Value = 1.5
$ws.Range("Y10").Value = 1.29
$ws.Range("AA10").Value = 1.53
$ws.Range("AB10").Value = 2.38
$ws.Range("AD10").Value = 11
$ws.Range("AF10").Value = 17
$ws.Range("AJ10").Value = 8
$ws.Range("G11").Value = 2.05
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 2.75
$ws.Range("M11").Value = 1.03
$ws.Range("N11").Value = 11
$ws.Range("O11").Value = 1.22
$ws.Range("S11").Value = 1.85
$ws.Range("T11").Value = 2
$ws.Range("X11").Value = 1.33
$ws.Range("AC11").Value = 8.5
$ws.Range("AE11").Value = 9
$ws.Range("AF11").Value = 19
$ws.Range("AG11").Value = 17
$ws.Range("AM11").Value = 151
$ws.Range("AP11").Value = 12
$ws.Range("AR11").Value = 26
$ws.Range("M12").Value = 17
$ws.Range("N12").Value = 1.03
$ws.Range("O12").Value = 1.11
$ws.Range("X12").Value = 1.58
$ws.Range("M13").Value = 1.02
$ws.Range("O13").Value = 1.13
$ws.Range("W13").Value = 2.37
$ws.Range("X13").Value = 1.5
$ws.Range("M14").Value = 1.03
$ws.Range("O14").Value = 1.22
$ws.Range("X14").Value = 1.33
$ws.Range("M15").Value = 1.01
$ws.Range("O15").Value = 1.11
$ws.Range("X15").Value = 1.63
$ws.Range("AF15").Value = 15
$ws.Range("AJ15").Value = 8.5
$ws.Range("AL15").Value = 34
$ws.Range("AR15").Value = 29
$ws.Range("M16").Value = 1.03
$ws.Range("O16").Value = 1.22
$ws.Range("X16").Value = 1.33
$ws.Range("G17").Value = 5.5
$ws.Range("H17").Value = 4.15
$ws.Range("I17").Value = 1.52
$ws.Range("J17").Value = 5.3
$ws.Range("K17").Value = 2.32
$ws.Range("L17").Value = 2.05
$ws.Range("O17").Value = 1.22
$ws.Range("P17").Value = 3.85
$ws.Range("T17").Value = 2.1
$ws.Range("W17").Value = 2.57
$ws.Range("Y17").Value = 1.33
$ws.Range("Z17").Value = 3.05
$ws.Range("AA17").Value = 1.75
$ws.Range("AB17").Value = 1.95
$ws.Range("G18").Value = 2.15
$ws.Range("I18").Value = 3.1
$ws.Range("J18").Value = 2.72
$ws.Range("K18").Value = 2.1
$ws.Range("L18").Value = 3.65
$ws.Range("O18").Value = 1.33
$ws.Range("P18").Value = 2.8
$ws.Range("S18").Value = 1.98
$ws.Range("T18").Value = 1.65
$ws.Range("W18").Value = 3.25
$ws.Range("X18").Value = 1.25
$ws.Range("Y18").Value = 1.4
$ws.Range("Z18").Value = 2.52
$ws.Range("AA18").Value = 1.8
$ws.Range("AB18").Value = 1.8
$ws.Range("AC18").Value = 7
$ws.Range("AD18").Value = 9.75
$ws.Range("AE18").Value = 9
$ws.Range("AF18").Value = 19.5
$ws.Range("AG18").Value = 18
$ws.Range("AI18").Value = 9
$ws.Range("AK18").Value = 15.5
$ws.Range("AL18").Value = 80
$ws.Range("AM18").Value = 700
$ws.Range("AN18").Value = 8.75
$ws.Range("AO18").Value = 15
$ws.Range("AP18").Value = 11.25
$ws.Range("AQ18").Value = 37
$ws.Range("AR18").Value = 28
$ws.Range("AS18").Value = 40
